$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 749, shifting everything
# (rows 749..790) down by two rows.
$ws.Range("A749:A750").EntireRow.Insert()

# Fill the two newly-inserted rows with the new data.
# Force column A to be stored as text (not auto-converted to a date
# serial number) the same way the rest of the date column is stored.
$dateRng = $ws.Range("A749:A750")
$dateRng.NumberFormat = "@"

$ws.Range("A749").Value = "2026/02/03"
$ws.Range("B749").Value = "火"
$ws.Range("C749").Value = 19
$ws.Range("D749").Value = 25

$ws.Range("A750").Value = "2026/02/03"
$ws.Range("B750").Value = "火"
$ws.Range("C750").Value = 22
$ws.Range("D750").Value = 24

# Drop the temporary number-format style so the cells fall back to the
# workbook's default (unstyled) cell, matching the rest of the sheet.
$dateRng.Style = "Normal"
